# Auto-generated edit script: refresh crypto price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.733.92'
$ws.Range('E2').Value = '  +1.54%  '

$ws.Range('D3').Value = '3.495.56'
$ws.Range('E3').Value = '  +0.33%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').Value = '''602.47'
$ws.Range('E5').Value = '  +1.90%  '

$ws.Range('D6').Value = '''172.58'
$ws.Range('E6').Value = '  +2.50%  '

$ws.Range('E7').Value = '  +0.10%  '

$ws.Range('D8').Value = '3.490.71'
$ws.Range('E8').Value = '  +0.30%  '

$ws.Range('E9').Value = '  +0.00%  '

$ws.Range('E10').Value = '  +0.09%  '

$ws.Range('D11').Value = '''7.29'
$ws.Range('E11').Value = '  +7.61%  '

$ws.Range('E12').Value = '  +1.21%  '

$ws.Range('D13').Value = '''45.95'
$ws.Range('E13').Value = '  -1.50%  '

$ws.Range('E14').Value = '  -1.75%  '

$ws.Range('D15').Value = '4.057.58'
$ws.Range('E15').Value = '  +0.47%  '

$ws.Range('E16').Value = '  -0.75%  '

$ws.Range('D17').Value = '''608.85'
$ws.Range('E17').Value = '  -0.88%  '

$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '69.724.39'
$ws.Range('E18').Value = '  +1.52%  '

$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.471.57'
$ws.Range('E19').Value = '  -0.35%  '

$ws.Range('E20').Value = '  +0.85%  '

$ws.Range('D21').Value = '''17.13'
$ws.Range('E21').Value = '  -0.53%  '

$ws.Range('D22').Value = '''0.869'
$ws.Range('E22').Value = '  -0.27%  '

$ws.Range('D23').Value = '''9.03'
$ws.Range('E23').Value = '  -18.99%  '

$ws.Range('D24').Value = '''15.43'
$ws.Range('E24').Value = '  -2.11%  '

$ws.Range('D25').Value = '''95.68'

$ws.Range('D26').Value = '''3.70'
$ws.Range('E26').Value = '  -2.12%  '

$ws.Range('D27').Value = '''1.00'
$ws.Range('E27').Value = '  -0.09%  '

$ws.Range('E28').Value = '  -2.28%  '

$ws.Range('D29').Value = '''34.23'
$ws.Range('E29').Value = '  +4.25%  '

$ws.Range('D30').Value = '''8.92'
$ws.Range('E30').Value = '  -2.21%  '

$ws.Range('D31').Value = '''658.51'
$ws.Range('E31').Value = '  +15.02%  '

$ws.Range('D32').Value = '''8.08'
$ws.Range('E32').Value = '  -3.95%  '

$ws.Range('D33').Value = '''2.95'
$ws.Range('E33').Value = '  -4.35%  '

$ws.Range('E34').Value = '  +1.31%  '

$ws.Range('E35').Value = '  -3.81%  '

$ws.Range('D36').Value = '''0.0995'
$ws.Range('E36').Value = '  -1.63%  '

$ws.Range('D37').Value = '''3.55'

$ws.Range('D38').Value = '''10.66'
$ws.Range('E38').Value = '  -0.36%  '

$ws.Range('D39').Value = '''0.0472'
$ws.Range('E39').Value = '  +8.10%  '

$ws.Range('E40').Value = '  +0.22%  '

$ws.Range('D41').Value = '''56.33'
$ws.Range('E41').Value = '  -1.16%  '

$ws.Range('D42').Value = '''0.142'
$ws.Range('E42').Value = '  +3.40%  '

$ws.Range('D43').Value = '3.309.05'
$ws.Range('E43').Value = '  -2.50%  '

$ws.Range('D44').Value = '''0.310'
$ws.Range('E44').Value = '  -4.07%  '

$ws.Range('E45').Value = '  +2.64%  '

$ws.Range('D46').Value = '''32.15'
$ws.Range('E46').Value = '  -1.22%  '

$ws.Range('D47').Value = '0.0₃0685'
$ws.Range('E47').Value = '  -0.79%  '

$ws.Range('E48').Value = '  -0.52%  '

$ws.Range('E49').Value = '  +0.83%  '

$ws.Range('D50').Value = '''134.08'
$ws.Range('E50').Value = '  +1.34%  '

